# Quarterly indexing esoteric bug-fix operation
# Updates the forecast error statistics (ME, MAE, MSE, RMSE, SE, N) for rows
# corresponding to quarters Q1..Q9 (sheet rows 2-11) to reflect the corrected
# quarterly indexing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1113813209086354, 1.082668180363185, 3.265531177099384, 1.807078077200701, 1.828521252897521, 37)
    3  = @(0.1513751093172216, 0.9545256880515617, 2.919874515134196, 1.7087640314374, 1.726189554985052, 36)
    4  = @(0.2122263332368638, 0.9540284968671798, 2.851160657205658, 1.688538023618555, 1.699603954074797, 35)
    5  = @(0.1511539997107805, 0.9530044890603636, 2.951427543481194, 1.717971927442702, 1.737044864330408, 34)
    6  = @(0.2330575124652057, 0.9893939392548374, 2.998246469334817, 1.731544532876593, 1.742391603086522, 33)
    7  = @(0.1708011111517192, 0.9843372303369592, 3.106494255942373, 1.762524966047963, 1.782298982592832, 32)
    8  = @(0.2456821311865818, 1.039250576912814, 3.172757163278162, 1.781223501775721, 1.793361100136716, 31)
    9  = @(0.1868069867954431, 1.037514270990646, 3.290947412145117, 1.814096858534603, 1.835300565959116, 30)
    10 = @(0.2405366504582205, 1.068640860714179, 3.366469277036136, 1.834794069381122, 1.851155339123688, 29)
    11 = @(0.1965726553058998, 1.065268685810892, 3.497083030370647, 1.870048937961423, 1.893814318830355, 28)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$row").Value = $values[$i]
    }
}
